$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decision table rule "R30" (row 10): the lower bound ("min") of the
# hour range for the "Good Evening" greeting is corrected from 18 to 1.
$ws.Range("C10").Value = 1
